$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing order row (row 2) with the latest values for this order
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "ser"
$ws.Range("C2").Value = "e"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 48.26

# Loop puts each new order on a new line - add the next order on row 3.
# Match the existing data rows' formatting (vertical-bottom alignment, same
# as row 2) instead of inheriting the blank column default formatting.
$ws.Range("B3:H3").VerticalAlignment = -4107

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "e"
$ws.Range("C3").Value = "aaaa"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 23
$ws.Range("H3").Value = 468.68
